# Updates odds values on Sheet1 to match the 2024-10-12 FlashScore refresh.
# Columns are addressed by (row, 1-based column index) via Cells.Item(row, col).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 15).Value = 1.29
$ws.Cells.Item(2, 16).Value = 3.5
$ws.Cells.Item(2, 17).Value = 1.98
$ws.Cells.Item(2, 18).Value = 1.88
$ws.Cells.Item(3, 7).Value = 1.65
$ws.Cells.Item(3, 9).Value = 6.25
$ws.Cells.Item(3, 10).Value = 2.3
$ws.Cells.Item(3, 17).Value = 2.3
$ws.Cells.Item(3, 18).Value = 1.6
$ws.Cells.Item(3, 21).Value = 2.25
$ws.Cells.Item(3, 22).Value = 1.57
$ws.Cells.Item(3, 28).Value = 41
$ws.Cells.Item(3, 30).Value = 7
$ws.Cells.Item(3, 35).Value = 21
$ws.Cells.Item(3, 40).Value = 3.4
$ws.Cells.Item(3, 47).Value = 10
$ws.Cells.Item(3, 52).Value = 151
$ws.Cells.Item(3, 53).Value = 201
$ws.Cells.Item(4, 13).Value = 1.14
$ws.Cells.Item(4, 14).Value = 5.5
$ws.Cells.Item(4, 15).Value = 1.62
$ws.Cells.Item(4, 16).Value = 2.2
$ws.Cells.Item(6, 10).Value = 3.6
$ws.Cells.Item(6, 11).Value = 1.83
$ws.Cells.Item(6, 13).Value = 1.14
$ws.Cells.Item(6, 14).Value = 5.5
$ws.Cells.Item(6, 15).Value = 1.57
$ws.Cells.Item(6, 16).Value = 2.25
$ws.Cells.Item(6, 17).Value = 2.88
$ws.Cells.Item(6, 18).Value = 1.4
$ws.Cells.Item(6, 21).Value = 2.25
$ws.Cells.Item(6, 22).Value = 1.57
$ws.Cells.Item(6, 25).Value = 12
$ws.Cells.Item(6, 27).Value = 29
$ws.Cells.Item(6, 29).Value = 5.5
$ws.Cells.Item(6, 34).Value = 13
$ws.Cells.Item(11, 7).Value = 4.1
$ws.Cells.Item(11, 8).Value = 3.15
$ws.Cells.Item(11, 11).Value = 1.98
$ws.Cells.Item(11, 18).Value = 1.52
$ws.Cells.Item(11, 19).Value = 1.47
$ws.Cells.Item(11, 20).Value = 2.35
$ws.Cells.Item(11, 23).Value = 9.25
$ws.Cells.Item(11, 24).Value = 21
$ws.Cells.Item(11, 27).Value = 50
$ws.Cells.Item(11, 29).Value = 7.2
$ws.Cells.Item(11, 31).Value = 18
$ws.Cells.Item(11, 34).Value = 8
$ws.Cells.Item(11, 38).Value = 37
$ws.Cells.Item(11, 41).Value = 25
$ws.Cells.Item(11, 42).Value = 35
$ws.Cells.Item(11, 47).Value = 7.7
$ws.Cells.Item(11, 50).Value = 9.5
$ws.Cells.Item(12, 10).Value = 2.3
$ws.Cells.Item(12, 12).Value = 4.2
$ws.Cells.Item(12, 17).Value = 1.6
$ws.Cells.Item(12, 25).Value = 8.25
$ws.Cells.Item(12, 27).Value = 12.5
$ws.Cells.Item(12, 33).Value = 14.5
$ws.Cells.Item(12, 34).Value = 26
$ws.Cells.Item(12, 39).Value = 300
$ws.Cells.Item(12, 42).Value = 16
$ws.Cells.Item(12, 44).Value = 55
$ws.Cells.Item(12, 47).Value = 6.8
$ws.Cells.Item(12, 49).Value = 5.9
$ws.Cells.Item(12, 50).Value = 21
$ws.Cells.Item(12, 51).Value = 24
$ws.Cells.Item(12, 54).Value = 250
$ws.Cells.Item(14, 7).Value = 2.45
$ws.Cells.Item(14, 8).Value = 2.9
$ws.Cells.Item(14, 9).Value = 2.88
$ws.Cells.Item(14, 10).Value = 3.4
$ws.Cells.Item(14, 12).Value = 3.75
$ws.Cells.Item(14, 13).Value = 1.11
$ws.Cells.Item(14, 14).Value = 6.5
$ws.Cells.Item(14, 15).Value = 1.53
$ws.Cells.Item(14, 16).Value = 2.38
$ws.Cells.Item(14, 17).Value = 2.7
$ws.Cells.Item(14, 18).Value = 1.44
$ws.Cells.Item(14, 19).Value = 1.62
$ws.Cells.Item(14, 20).Value = 2.2
$ws.Cells.Item(14, 24).Value = 11
$ws.Cells.Item(14, 25).Value = 11
$ws.Cells.Item(14, 26).Value = 26
$ws.Cells.Item(14, 27).Value = 26
$ws.Cells.Item(14, 29).Value = 6
$ws.Cells.Item(14, 42).Value = 34
$ws.Cells.Item(14, 44).Value = 101
$ws.Cells.Item(14, 46).Value = 2.2
$ws.Cells.Item(14, 47).Value = 9.5
$ws.Cells.Item(14, 49).Value = 4.75
$ws.Cells.Item(20, 10).Value = 2.4
$ws.Cells.Item(20, 12).Value = 5.5
$ws.Cells.Item(20, 13).Value = 1.07
$ws.Cells.Item(20, 14).Value = 9
$ws.Cells.Item(20, 26).Value = 13
$ws.Cells.Item(20, 34).Value = 23
$ws.Cells.Item(20, 41).Value = 9.5
$ws.Cells.Item(20, 43).Value = 34
$ws.Cells.Item(20, 52).Value = 101
$ws.Cells.Item(20, 54).Value = 351
$ws.Cells.Item(21, 17).Value = 2.3
$ws.Cells.Item(21, 18).Value = 1.6
$ws.Cells.Item(22, 7).Value = 3.5
$ws.Cells.Item(22, 8).Value = 3.25
$ws.Cells.Item(22, 9).Value = 2.15
$ws.Cells.Item(22, 10).Value = 4
$ws.Cells.Item(22, 12).Value = 2.88
$ws.Cells.Item(22, 15).Value = 1.4
$ws.Cells.Item(22, 16).Value = 2.75
$ws.Cells.Item(22, 17).Value = 2.25
$ws.Cells.Item(22, 18).Value = 1.62
$ws.Cells.Item(22, 21).Value = 1.91
$ws.Cells.Item(22, 22).Value = 1.8
$ws.Cells.Item(22, 24).Value = 17
$ws.Cells.Item(22, 25).Value = 13
$ws.Cells.Item(22, 26).Value = 41
$ws.Cells.Item(22, 29).Value = 8.5
$ws.Cells.Item(22, 30).Value = 6.5
$ws.Cells.Item(22, 31).Value = 17
$ws.Cells.Item(22, 34).Value = 9.5
$ws.Cells.Item(22, 36).Value = 19
$ws.Cells.Item(22, 37).Value = 19
$ws.Cells.Item(22, 39).Value = 351
$ws.Cells.Item(22, 40).Value = 5.5
$ws.Cells.Item(22, 41).Value = 21
$ws.Cells.Item(22, 43).Value = 67
$ws.Cells.Item(22, 44).Value = 101
$ws.Cells.Item(22, 45).Value = 251
$ws.Cells.Item(22, 49).Value = 4
$ws.Cells.Item(22, 50).Value = 12
$ws.Cells.Item(25, 13).Value = 1.05
$ws.Cells.Item(25, 14).Value = 11
$ws.Cells.Item(25, 15).Value = 1.29
$ws.Cells.Item(25, 16).Value = 3.5
$ws.Cells.Item(25, 17).Value = 1.98
$ws.Cells.Item(25, 18).Value = 1.88
$ws.Cells.Item(29, 8).Value = 3.55
$ws.Cells.Item(29, 9).Value = 2.5
$ws.Cells.Item(29, 10).Value = 3
$ws.Cells.Item(29, 30).Value = 7
$ws.Cells.Item(29, 32).Value = 55
$ws.Cells.Item(29, 36).Value = 27
$ws.Cells.Item(29, 37).Value = 19.5
$ws.Cells.Item(29, 44).Value = 80
$ws.Cells.Item(29, 47).Value = 7
$ws.Cells.Item(29, 50).Value = 13
$ws.Cells.Item(31, 10).Value = 2.25
$ws.Cells.Item(31, 11).Value = 2.2
$ws.Cells.Item(31, 19).Value = 1.38
$ws.Cells.Item(31, 20).Value = 2.8
$ws.Cells.Item(31, 27).Value = 13
$ws.Cells.Item(31, 41).Value = 8.25
$ws.Cells.Item(31, 42).Value = 17.5
$ws.Cells.Item(31, 46).Value = 2.8
$ws.Cells.Item(37, 8).Value = 2.9
$ws.Cells.Item(37, 9).Value = 3.15
$ws.Cells.Item(37, 10).Value = 2.9
$ws.Cells.Item(37, 11).Value = 1.98
$ws.Cells.Item(37, 15).Value = 1.37
$ws.Cells.Item(37, 16).Value = 2.65
$ws.Cells.Item(37, 17).Value = 2.07
$ws.Cells.Item(37, 18).Value = 1.6
$ws.Cells.Item(37, 19).Value = 1.42
$ws.Cells.Item(37, 20).Value = 2.45
$ws.Cells.Item(37, 21).Value = 1.78
$ws.Cells.Item(37, 22).Value = 1.83
$ws.Cells.Item(37, 23).Value = 7.1
$ws.Cells.Item(37, 24).Value = 11.25
$ws.Cells.Item(37, 25).Value = 9
$ws.Cells.Item(37, 27).Value = 20
$ws.Cells.Item(37, 28).Value = 32
$ws.Cells.Item(37, 29).Value = 7.6
$ws.Cells.Item(37, 30).Value = 5.7
$ws.Cells.Item(37, 32).Value = 70
$ws.Cells.Item(37, 33).Value = 8.5
$ws.Cells.Item(37, 34).Value = 16
$ws.Cells.Item(37, 38).Value = 40
$ws.Cells.Item(37, 41).Value = 12
$ws.Cells.Item(37, 42).Value = 19.5
$ws.Cells.Item(37, 43).Value = 50
$ws.Cells.Item(37, 44).Value = 80
$ws.Cells.Item(37, 45).Value = 250
$ws.Cells.Item(37, 46).Value = 2.42
$ws.Cells.Item(37, 47).Value = 6.7
$ws.Cells.Item(37, 50).Value = 18
